$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at row 223, pushing existing rows 223-302 down to 225-304.
$ws.Range("223:224").Insert()

# Populate new row 223.
$ws.Cells.Item(223, 1).Value = 8
$ws.Cells.Item(223, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(223, 3).Value = "Coquimbo"
$ws.Cells.Item(223, 4).Value = 44900
$ws.Cells.Item(223, 5).Value = 4
$ws.Cells.Item(223, 6).Value = 100112031
$ws.Cells.Item(223, 7).Value = "Poroto verde"
$ws.Cells.Item(223, 8).Value = "Magnum"
$ws.Cells.Item(223, 9).Value = "Primera"
$ws.Cells.Item(223, 10).Value = 520
$ws.Cells.Item(223, 11).Value = 31000
$ws.Cells.Item(223, 12).Value = 32000
$ws.Cells.Item(223, 13).Value = 31500
$ws.Cells.Item(223, 14).Value = "`$/caja 25 kilos"
$ws.Cells.Item(223, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(223, 16).Value = 1260
$ws.Cells.Item(223, 17).Value = 25
$ws.Cells.Item(223, 18).Value = "Hortaliza"

# Populate new row 224.
$ws.Cells.Item(224, 1).Value = 8
$ws.Cells.Item(224, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(224, 3).Value = "Coquimbo"
$ws.Cells.Item(224, 4).Value = 44900
$ws.Cells.Item(224, 5).Value = 4
$ws.Cells.Item(224, 6).Value = 100112031
$ws.Cells.Item(224, 7).Value = "Poroto verde"
$ws.Cells.Item(224, 8).Value = "Magnum"
$ws.Cells.Item(224, 9).Value = "Primera"
$ws.Cells.Item(224, 10).Value = 400
$ws.Cells.Item(224, 11).Value = 25000
$ws.Cells.Item(224, 12).Value = 26000
$ws.Cells.Item(224, 13).Value = 25500
$ws.Cells.Item(224, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(224, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(224, 16).Value = 1020
$ws.Cells.Item(224, 17).Value = 25
$ws.Cells.Item(224, 18).Value = "Hortaliza"

Write-Host "done"
